# Phase1-Backup-Guide.pptx edit
#
# 1. Slide 7 ("Step 5: Outlook Data Files (OST/PST)"):
#    - Reposition/resize the red warning textbox (TextBox 2)
#    - Append "!" to the warning text
#    - Move the screenshot picture down slightly
# 2. Delete slide 8 ("Step 6: iOS Backup Verification") entirely -
#    "Phase 1 Checklist" becomes the new (last) slide.

$p = $ppt.ActivePresentation

# --- Slide 7: Step 5 (Outlook Data Files) tweaks ---
$s7 = $p.Slides.Item(7)

# Warning textbox: reposition + resize
$warnBox = $s7.Shapes.Item(2)
$warnBox.Left   = 156.5827
$warnBox.Top    = 129.6
$warnBox.Width  = 406.8346
$warnBox.Height = 26.6578

# Append "!" to the warning text (creates a second run, matching the edit)
$null = $warnBox.TextFrame.TextRange.InsertAfter("!")

# Screenshot picture: move down
$pic = $s7.Shapes.Item(4)
$pic.Top = 151.20004

# --- Remove Slide 8: Step 6 iOS Backup Verification ---
$s8 = $p.Slides.Item(8)
$s8.Delete()
